# -----------------------------------------------------------------------
# Unicode helpers (literal non-ASCII chars get mangled by the PS host's
# script loader, so build them from code points instead).
# -----------------------------------------------------------------------
$ENDASH  = [char]0x2013   # "-" (en dash)
$EMDASH  = [char]0x2014   # "-" (em dash)
$LDQUOTE = [char]0x201C   # left double quotation mark
$RDQUOTE = [char]0x201D   # right double quotation mark
$OSLASH  = [char]0x00D8   # O with stroke
$OUML    = [char]0x00F6   # o with diaeresis

$d = $word.ActiveDocument

# -----------------------------------------------------------------------
# 1) Intro paragraph: expand "In this paper however we mainly focus only
#    on the Silicon Carbide MOSFET." into the longer, reworded passage
#    (several extra sourced sentences), then drop the paragraph that used
#    to carry the "These emerging ... [1]." sentence plus the blank
#    paragraph that followed it (that sentence now lives inside the
#    paragraph we just rewrote, renumbered to reference [3]).
# -----------------------------------------------------------------------
$rng = $d.Content
$rng.Find.Execute("In this paper however we mainly focus only on the Silicon Carbide MOSFET. ", `
    $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $rng.Find.Found) { throw "intro sentence not found" }

$newIntro = "In this paper however, we mainly focus only on the Silicon Carbide based Power devices. " + `
  "There has been a tremendous amount research effort on developing power semiconductor devices with Silicon Carbide (SiC) in the pursuit of higher efficiency and smaller dimensions [1], [2]. " + `
  "The availability of SiC wafers on a commercial basis has led to the demonstration of many types of metal-oxide semiconductor (MOS)-gated devices that exploit its unique properties. " + `
  "These emerging Silicon Carbide (SiC) MOSFET power devices promise to displace Silicon IGBTs from the majority of challenging power electronics applications by enabling superior efficiency and power density, as well as capability to operate at higher temperatures [3]. " + `
  "Reference [4] focuses on the comparison of a SiC based DC/DC converter and an IGBT based DC/DC converter and thus concludes that the efficiency of an SiC converter is greater than that of the IGBT converter over an output power range. " + `
  "An electro-thermal analysis of an automotive traction inverter platform based on SiC MOSFET and SiC IGBT technology is discussed in [5] and the results show that there is a higher total loss reduction in the SiC MOSFET model compared to the IGBT model. "
$rng.Text = $newIntro

# Now remove the old "These emerging ... [1]." paragraph and the blank
# paragraph right after it (the blank paragraph that remains right after
# the rewritten intro paragraph is left untouched).
$rng2 = $d.Content
$rng2.Find.Execute("These emerging Silicon Carbide (SiC) MOSFET power devices promise to displace Silicon IGBTs from the majority of challenging power electronics applications by enabling superior efficiency and power density, as well as capability to operate at higher temperatures [1]. ", `
    $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $rng2.Find.Found) { throw "old 'these emerging' sentence not found" }
$oldPara = $rng2.Paragraphs(1)
$blankPara = $oldPara.Next()
$killRange = $d.Range($oldPara.Range.Start, $blankPara.Range.End)
$killRange.Delete()

# -----------------------------------------------------------------------
# 2) Cycloconverter paragraph: merge the two runs that used to be split
#    by a lastRenderedPageBreak into one contiguous run (no text change).
# -----------------------------------------------------------------------
$rng3 = $d.Content
$rng3.Find.Execute(" A Cycloconverter is a device that converts constant voltage and frequency AC waveform to another AC waveform of lower frequency without using DC link in the conversion process th", `
    $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $rng3.Find.Found) {
    # Not merged yet (original has the page break splitting "DC " / "link ...th")
    $rng3b = $d.Content
    $rng3b.Find.Execute(" A Cycloconverter is a device that converts constant voltage and frequency AC waveform to another AC waveform of lower frequency without using DC link in the conversion process th", `
        $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
}
$rngMerge = $d.Content
$found = $rngMerge.Find.Execute(" A Cycloconverter is a device that converts constant voltage and frequency AC waveform to another AC waveform of lower frequency without using DC *link in the conversion process th", `
    $false, $false, $false, $false, $false, $true, 1, $true, "", 0)
if ($found) {
    $rngMerge.Text = " A Cycloconverter is a device that converts constant voltage and frequency AC waveform to another AC waveform of lower frequency without using DC link in the conversion process th"
}

# -----------------------------------------------------------------------
# 3) "Figure Labels" paragraph: merge the two runs split by a
#    lastRenderedPageBreak into a single run (no text change).
# -----------------------------------------------------------------------
$rng4 = $d.Content
$rng4.Find.Execute("Figure Labels: Use 8 point Times New Roman for Figure labels. Use words rather than symbols or abbreviations when writing Figure axis labels to avoid *reader. As an example*Temperature/K" + [char]0x201D + ".", `
    $false, $false, $false, $false, $false, $true, 1, $true, "", 0)
if ($rng4.Find.Found) {
    $mag  = $LDQUOTE + "Magnetization" + $RDQUOTE
    $mag2 = $LDQUOTE + "Magnetization, M" + $RDQUOTE
    $mJ   = $LDQUOTE + "M" + $RDQUOTE
    $magAm  = $LDQUOTE + "Magnetization (A/m)" + $RDQUOTE
    $magBr  = $LDQUOTE + "Magnetization {A[m(1)]}" + $RDQUOTE
    $am     = $LDQUOTE + "A/m" + $RDQUOTE
    $tempK  = $LDQUOTE + "Temperature (K)" + $RDQUOTE
    $tempOverK = $LDQUOTE + "Temperature/K" + $RDQUOTE
    $figTxt = "Figure Labels: Use 8 point Times New Roman for Figure labels. Use words rather than symbols or abbreviations when writing Figure axis labels to avoid confusing the reader. As an example, write the quantity $mag, or $mag2, not just $mJ. If including units in the label, present them within parentheses. Do not label axes only with units. In the example, write $magAm or $magBr, not just $am. Do not label axes with a ratio of quantities and units. For example, write $tempK, not $tempOverK."
    $rng4.Text = $figTxt
}

# -----------------------------------------------------------------------
# 4) References list rework.
# -----------------------------------------------------------------------

# 4a) Insert a brand-new reference (McBryde et al.) right before the
#     "L. D. Stevanovic ..." reference paragraph.
$rngStev = $d.Content
$rngStev.Find.Execute("L. D. Stevanovic, K. S. Matocha, P. A. Losee, J. S. Glaser, J. J. Nasadoski and S. D. Arthur", `
    $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $rngStev.Find.Found) { throw "Stevanovic reference not found" }
$stevPara = $rngStev.Paragraphs(1)
$insertPoint = $stevPara.Range.Duplicate
$insertPoint.Collapse(1)
$insertPoint.InsertParagraphBefore()
$newRefPara = $insertPoint.Paragraphs(1)
$newRefRange = $newRefPara.Range
$newRefRange.Text = "J. McBryde, A. Kadavelugu, B. Compton, S. Bhattacharya, M. Das and A. Agarwal, " + $LDQUOTE + "Performance comparison of 1200V Silicon and SiC devices for UPS application," + $RDQUOTE + " IECON 2010 - 36th Annual Conference on IEEE Industrial Electronics Society, Glendale, AZ, 2010, pp. 2657-2662."
$newRefPara.Range.ParagraphStyle = "List Paragraph"
$newRefPara.Range.ParagraphFormat.CharacterUnitLeftIndent = 0
$ilvl = $newRefPara.Range.ListFormat

# Apply the numbered-list formatting (ilvl 0 / numId 1) to match the rest
# of the reference list, then apply Emphasis to the venue-name substring.
$newRefPara.Range.ListFormat.ApplyNumberDefault()

$rngEmph = $d.Content
$rngEmph.Find.Execute("IECON 2010 - 36th Annual Conference on IEEE Industrial Electronics Society", `
    $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($rngEmph.Find.Found) {
    $rngEmph.Style = "Emphasis"
}

# 4b) Replace the (now second) "L. D. Stevanovic..." reference's text with
#     the Biela et al. citation.
$rngStev2 = $d.Content
$rngStev2.Find.Execute("L. D. Stevanovic, K. S. Matocha, P. A. Losee, J. S. Glaser, J. J. Nasadoski and S. D. Arthur, " + [char]0x22 + "Recent advances in silicon carbide MOSFET power devices," + [char]0x22 + " 2010 Twenty-Fifth Annual IEEE Applied Power Electronics Conference and Exposition (APEC), Palm Springs, CA, 2010, pp. 401-407.", `
    $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $rngStev2.Find.Found) { throw "Stevanovic full reference text not found" }
$bielaText = "J. Biela, M. Schweizer, S. Waffler and J. W. Kolar, " + $LDQUOTE + "SiC versus Si" + $EMDASH + "Evaluation of Potentials for Performance Improvement of Inverter and DC" + $ENDASH + "DC Converter Systems by SiC Power Semiconductors," + $RDQUOTE + " in IEEE Transactions on Industrial Electronics, vol. 58, no. 7, pp. 2872-2882, July 2011."
$rngStev2.Text = $bielaText
$rngEmph2 = $d.Content
$rngEmph2.Find.Execute("IEEE Transactions on Industrial Electronics", `
    $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($rngEmph2.Find.Found) {
    $rngEmph2.Style = "Emphasis"
}

# 4c) Replace "J. Clerk Maxwell..." reference with the original Stevanovic
#     citation (re-inserted here, after the Biela reference).
$rngMax = $d.Content
$rngMax.Find.Execute("J. Clerk Maxwell, A Treatise on Electricity and Magnetism, 3rd ed., vol. 2. Oxford: Clarendon, 1892, pp.68" + $ENDASH + "73.", `
    $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $rngMax.Find.Found) { throw "Maxwell reference not found" }
$stevText = "L. D. Stevanovic, K. S. Matocha, P. A. Losee, J. S. Glaser, J. J. Nasadoski and S. D. Arthur, " + $LDQUOTE + "Recent advances in silicon carbide MOSFET power devices," + $RDQUOTE + " 2010 Twenty-Fifth Annual IEEE Applied Power Electronics Conference and Exposition (APEC), Palm Springs, CA, 2010, pp. 401-407."
$rngMax.Text = $stevText
$rngMax.Font.Reset()
$rngEmph3 = $d.Content
$rngEmph3.Find.Execute("2010 Twenty-Fifth Annual IEEE Applied Power Electronics Conference and Exposition (APEC)", `
    $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($rngEmph3.Find.Found) {
    $rngEmph3.Style = "Emphasis"
}

# 4d) Replace "I. S. Jacobs..." reference with the Nielsen et al. citation.
$rngJac = $d.Content
$rngJac.Find.Execute("I. S. Jacobs and C. P. Bean, " + $LDQUOTE + "Fine particles, thin films and exchange anisotropy," + $RDQUOTE + " in Magnetism, vol. III, G. T. Rado and H. Suhl, Eds. New York: Academic, 1963, pp. 271" + $ENDASH + "350.", `
    $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $rngJac.Find.Found) { throw "Jacobs reference not found" }
$nielsenText = "R. " + $OSLASH + ". Nielsen, L. T" + $OUML + "r" + $OUML + "k, S. Munk-Nielsen and F. Blaabjerg, " + $LDQUOTE + "Efficiency and cost comparison of Si IGBT and SiC JFET isolated DC/DC converters," + $RDQUOTE + " IECON 2013 - 39th Annual Conference of the IEEE Industrial Electronics Society, Vienna, 2013, pp. 695-699."
$rngJac.Text = $nielsenText
$rngJac.Font.Reset()
$rngEmph4 = $d.Content
$rngEmph4.Find.Execute("IECON 2013 - 39th Annual Conference of the IEEE Industrial Electronics Society", `
    $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($rngEmph4.Find.Found) {
    $rngEmph4.Style = "Emphasis"
}

# 4e) Replace "K. Elissa..." reference with the Kempitiya & Chou citation.
$rngEli = $d.Content
$rngEli.Find.Execute("K. Elissa, " + $LDQUOTE + "Title of paper if known," + $RDQUOTE + " unpublished.", `
    $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $rngEli.Find.Found) { throw "Elissa reference not found" }
$kempText = "A. Kempitiya and W. Chou, " + $LDQUOTE + "An electro-thermal performance analysis of SiC MOSFET vs Si IGBT and diode automotive traction inverters under various drive cycles," + $RDQUOTE + " 2018 34th Thermal Measurement, Modeling & Management Symposium (SEMI-THERM), San Jose, CA, 2018, pp. 213-217."
$rngEli.Text = $kempText
$rngEli.Font.Reset()
$rngEmph5 = $d.Content
$rngEmph5.Find.Execute("2018 34th Thermal Measurement, Modeling & Management Symposium (SEMI-THERM)", `
    $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($rngEmph5.Find.Found) {
    $rngEmph5.Style = "Emphasis"
}

# 4f) Empty out the "R. Nicole..." reference: drop its text and its list
#     numbering (it becomes a bare, unnumbered paragraph).
$rngNic = $d.Content
$rngNic.Find.Execute("R. Nicole, " + $LDQUOTE + "Title of paper with only first word capitalized," + $RDQUOTE + " J. Name Stand. Abbrev., in press.", `
    $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $rngNic.Find.Found) { throw "Nicole reference not found" }
$nicPara = $rngNic.Paragraphs(1)
$nicPara.Range.Text = ""
$nicPara.Range.ListFormat.RemoveNumbers()
$nicPara.LeftIndent = 354
$nicPara.FirstLineIndent = 0

# 4g) Delete the "Y. Yorozu..." and "M. Young..." references entirely
#     (both paragraphs, including their paragraph marks).
$rngYor = $d.Content
$rngYor.Find.Execute("Y. Yorozu, M. Hirano, K. Oka, and Y. Tagawa", `
    $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $rngYor.Find.Found) { throw "Yorozu reference not found" }
$yorPara = $rngYor.Paragraphs(1)
$youngPara = $yorPara.Next()
$killRange2 = $d.Range($yorPara.Range.Start, $youngPara.Range.End)
$killRange2.Delete()
